$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D:K -> E:L), matching a newly
# reported quarter (period ending 2018-09-30) being added to the report.
$ws.Columns("D").Insert()

# Carry the per-row number formatting from the (now-shifted) column E into
# the freshly inserted, blank column D.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column with the latest quarter's figures.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 61500
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -3100
$ws.Range("D17").Value = 20900
$ws.Range("D18").Value = 40600
$ws.Range("D20").Value = -32200
$ws.Range("D21").Value = 11500
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 8400
$ws.Range("D24").Value = -1000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 9400
$ws.Range("D27").Value = 9000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 32200
$ws.Range("D33").Value = 9000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 9000
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 57100
$ws.Range("D42").Value = 493400
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 25000
$ws.Range("D49").Value = 372200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 22800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7395100
$ws.Range("D57").Value = 2400
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 132900
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 6358000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 29100
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 271300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1007900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 9000
$ws.Range("D83").Value = 3000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 43300
$ws.Range("D91").Value = -900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -114900
$ws.Range("D96").Value = -4100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 203800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 132300

Write-Output "done"
